$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values must stay text - force Text format before assigning
# so Excel COM does not auto-coerce strings like "583.78" into numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.171.08"
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.487.22"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.78"
$ws.Range("E5").Value = "  -0.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.22"
$ws.Range("E6").Value = "  +2.89%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  -0.68%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.486.90"
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("E10").Value = "  +0.70%  "
$ws.Range("E11").Value = "  -0.05%  "
$ws.Range("E12").Value = "  -0.50%  "
$ws.Range("E13").Value = "  -2.24%  "
$ws.Range("E14").Value = "  +1.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.33"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.091.78"
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("E17").Value = "  -1.82%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.539.96"
$ws.Range("E18").Value = "  +2.04%  "
$ws.Range("E19").Value = "  -5.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.39"
$ws.Range("E20").Value = "  -5.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "348.59"
$ws.Range("E21").Value = "  -3.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.04"
$ws.Range("E22").Value = "  -1.97%  "
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("E24").Value = "  -4.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "68.43"
$ws.Range("E25").Value = "  -2.96%  "
$ws.Range("E26").Value = "  -2.94%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.28"
$ws.Range("E27").Value = "  -1.58%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("E30").Value = "  -3.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "510.99"
$ws.Range("E31").Value = "  +2.37%  "
$ws.Range("E32").Value = "  -3.95%  "
$ws.Range("E33").Value = "  -2.74%  "
$ws.Range("E34").Value = "  -3.79%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "160.38"
$ws.Range("E36").Value = "  +0.66%  "
$ws.Range("E37").Value = "  -7.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.71"
$ws.Range("E38").Value = "  +0.81%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.25"
$ws.Range("E39").Value = "  -4.05%  "
$ws.Range("E40").Value = "  -5.51%  "
$ws.Range("E41").Value = "  -2.57%  "
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.329"
$ws.Range("E43").Value = "  -2.16%  "
$ws.Range("E44").Value = "  -2.66%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.37"
$ws.Range("E45").Value = "  -4.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "38.80"
$ws.Range("E46").Value = "  -1.38%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "142.77"
$ws.Range("E47").Value = "  +0.70%  "
$ws.Range("E48").Value = "  -4.33%  "
$ws.Range("E49").Value = "  -4.76%  "
$ws.Range("E50").Value = "  -4.78%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0730"
$ws.Range("E51").Value = "  -1.02%  "

Write-Host "Applied cryptos list update"
